$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: add the "To" time and the computed "Time" formula
$ws.Range("C38").Value = 0.42708333333333331
$ws.Range("D38").Formula = "=C38-B38"

# Row 39: new entry "ROM 1" testbench work on 4.4.2020
$ws.Range("A39").Value = "4.4.2020"
$ws.Range("B39").Value = 0.4375
$ws.Range("C39").Value = 0.44791666666666669
$ws.Range("D39").Formula = "=C39-B39"
$ws.Range("E39").Value = "ROM 1"
$ws.Range("F39").Value = "Testbench"

# Match the styles of the row above (row 38) for the new row 39 cells
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B38:C38").Copy()
$ws.Range("B39:C39").PasteSpecial(-4122)

$ws.Range("D37").Copy()
$ws.Range("D39").PasteSpecial(-4122)

$ws.Range("E38:F38").Copy()
$ws.Range("E39:F39").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to match the authored state (G39 selected)
$ws.Range("G39").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
